$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 680 (the "ここに見えるのは何？" / ペンギン post) entirely;
# all subsequent rows shift up by one.
$ws.Rows.Item(680).Delete()
